# NIT-9006145760.xlsx -- Hoja1 "Periodo Mora" table refresh
# Reorders/interleaves periods 1706-2009 for both workers, adding
# PATRICIA MARRUGO IBARRA rows for periods 1706-1804 (part 1 of new EC data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    ,@(16, "1128044696", "ROSALBA TAJAN CERVANTES", "1706", 27578, 781242)
    ,@(17, "45507367", "PATRICIA MARRUGO IBARRA", "1706", 160000, 4000000)
    ,@(18, "1128044696", "ROSALBA TAJAN CERVANTES", "1707", 27578, 781242)
    ,@(19, "45507367", "PATRICIA MARRUGO IBARRA", "1707", 160000, 4000000)
    ,@(20, "1128044696", "ROSALBA TAJAN CERVANTES", "1708", 29509, 781242)
    ,@(21, "45507367", "PATRICIA MARRUGO IBARRA", "1708", 160000, 4000000)
    ,@(22, "1128044696", "ROSALBA TAJAN CERVANTES", "1709", 29509, 781242)
    ,@(23, "45507367", "PATRICIA MARRUGO IBARRA", "1709", 160000, 4000000)
    ,@(24, "1128044696", "ROSALBA TAJAN CERVANTES", "1710", 29509, 781242)
    ,@(25, "45507367", "PATRICIA MARRUGO IBARRA", "1710", 160000, 4000000)
    ,@(26, "1128044696", "ROSALBA TAJAN CERVANTES", "1711", 29509, 781242)
    ,@(27, "45507367", "PATRICIA MARRUGO IBARRA", "1711", 160000, 4000000)
    ,@(28, "1128044696", "ROSALBA TAJAN CERVANTES", "1712", 29509, 781242)
    ,@(29, "45507367", "PATRICIA MARRUGO IBARRA", "1712", 160000, 4000000)
    ,@(30, "1128044696", "ROSALBA TAJAN CERVANTES", "1801", 29509, 781242)
    ,@(31, "45507367", "PATRICIA MARRUGO IBARRA", "1801", 160000, 4000000)
    ,@(32, "1128044696", "ROSALBA TAJAN CERVANTES", "1802", 29509, 781242)
    ,@(33, "45507367", "PATRICIA MARRUGO IBARRA", "1802", 160000, 4000000)
    ,@(34, "1128044696", "ROSALBA TAJAN CERVANTES", "1803", 29509, 781242)
    ,@(35, "45507367", "PATRICIA MARRUGO IBARRA", "1803", 160000, 4000000)
    ,@(36, "1128044696", "ROSALBA TAJAN CERVANTES", "1804", 29509, 781242)
    ,@(37, "45507367", "PATRICIA MARRUGO IBARRA", "1804", 160000, 4000000)
    ,@(38, "1128044696", "ROSALBA TAJAN CERVANTES", "1901", 31249, 781242)
    ,@(39, "1128044696", "ROSALBA TAJAN CERVANTES", "1902", 31249, 781242)
    ,@(40, "1128044696", "ROSALBA TAJAN CERVANTES", "1903", 31249, 781242)
    ,@(41, "1128044696", "ROSALBA TAJAN CERVANTES", "1904", 31249, 781242)
    ,@(42, "1128044696", "ROSALBA TAJAN CERVANTES", "1905", 31249, 781242)
    ,@(43, "1128044696", "ROSALBA TAJAN CERVANTES", "1906", 31249, 781242)
    ,@(44, "1128044696", "ROSALBA TAJAN CERVANTES", "1907", 31249, 781242)
    ,@(45, "1128044696", "ROSALBA TAJAN CERVANTES", "1908", 31249, 781242)
    ,@(46, "1128044696", "ROSALBA TAJAN CERVANTES", "1909", 31249, 781242)
    ,@(47, "1128044696", "ROSALBA TAJAN CERVANTES", "1910", 31249, 781242)
    ,@(48, "1128044696", "ROSALBA TAJAN CERVANTES", "1911", 31249, 781242)
    ,@(49, "1128044696", "ROSALBA TAJAN CERVANTES", "1912", 31249, 781242)
    ,@(50, "1128044696", "ROSALBA TAJAN CERVANTES", "2001", 31249, 781242)
    ,@(51, "1128044696", "ROSALBA TAJAN CERVANTES", "2002", 31249, 781242)
    ,@(52, "1128044696", "ROSALBA TAJAN CERVANTES", "2003", 31249, 781242)
    ,@(53, "1128044696", "ROSALBA TAJAN CERVANTES", "2004", 31249, 781242)
    ,@(54, "1128044696", "ROSALBA TAJAN CERVANTES", "2005", 31249, 781242)
    ,@(55, "1128044696", "ROSALBA TAJAN CERVANTES", "2006", 31249, 781242)
    ,@(56, "1128044696", "ROSALBA TAJAN CERVANTES", "2007", 31249, 781242)
    ,@(57, "1128044696", "ROSALBA TAJAN CERVANTES", "2008", 31249, 781242)
    ,@(58, "1128044696", "ROSALBA TAJAN CERVANTES", "2009", 29166, 781242)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
